$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit swaps the position of the "e39f574c-..." and "72446337-..." file
# rows (rows 8 & 9) on all three sheets (Overview, zh-cn, de-de), and bumps
# the handoff/handback timestamp for the "e39f574c-..." row to reflect that
# it is now (also) "Ready for handoff" as of a later run.
# ---------------------------------------------------------------------------

# =====================  Sheet "Overview"  ===================================
$ws = $wb.Worksheets.Item("Overview")

# Row 8 becomes the 72446337 file (previously row 9), already "Ready for handoff"
$ws.Range("A8").Value = "72446337-f112-42e2-8cfc-c5b42674408a.md"
$ws.Range("B8").Value = "e2e\72446337-f112-42e2-8cfc-c5b42674408a.md"
$ws.Range("C8").Value = ".md"
$ws.Range("E8").Value = "Ready for handoff"
$ws.Range("F8").Value = "Ready for handoff"
$ws.Range("G8").Value = "2016-09-07 02:55:41"

# Row 9 becomes the e39f574c file (previously row 8), now also "Ready for handoff"
# with an updated (later) generate date
$ws.Range("A9").Value = "e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$ws.Range("B9").Value = "e2e\e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$ws.Range("C9").Value = ".md"
$ws.Range("E9").Value = "Ready for handoff"
$ws.Range("F9").Value = "Ready for handoff"
$ws.Range("G9").Value = "2016-09-07 03:08:41"

# Rebuild the hyperlinks on this sheet so that B8/B9 display text matches the
# (now swapped) file names while keeping the same underlying targets.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe85018917a4b35c6d742d73d641a78584ef05a/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/063ede46-8130-4a50-9e03-6494f895b9aa.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\063ede46-8130-4a50-9e03-6494f895b9aa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6e783c168c920fc8ee2c92b570627276ee1d6b/e2e/455c1867-16d2-41d4-b97b-ebfc6137f378.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\455c1867-16d2-41d4-b97b-ebfc6137f378.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/b9d5e119-8939-4ab9-ac75-292046d3d8d1.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\b9d5e119-8939-4ab9-ac75-292046d3d8d1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/897158552f04e39013c8f626cae66f5d9e2f0a5b/e2e/dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d3feedbe1faaf56c3e3a45b046ea53ccb6f815d/e2e/df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d482a0e3156d40b34d92c6dcf3f49c866658d47/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\72446337-f112-42e2-8cfc-c5b42674408a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49c8ff281eaf9230f89801459ee5727714e6c121/e2e/72446337-f112-42e2-8cfc-c5b42674408a.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\e39f574c-b915-4285-95c4-dfdd9da38f93.md") | Out-Null

# =====================  Sheet "zh-cn"  ======================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A8").Value = "72446337-f112-42e2-8cfc-c5b42674408a.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "72446337-f112-42e2-8cfc-c5b42674408a.24fd16ff902a084a4561b67b78c898624ad87e36.zh-cn.xlf"
$ws.Range("H8").Value = "2016-09-07 02:55:30"

$ws.Range("A9").Value = "e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("G9").Value = "e39f574c-b915-4285-95c4-dfdd9da38f93.6f9994cd4a132697e5c945d3f283e9462aaed486.zh-cn.xlf"
$ws.Range("H9").Value = "2016-09-07 03:08:29"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe85018917a4b35c6d742d73d641a78584ef05a/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", [System.Type]::Missing, [System.Type]::Missing, "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/736c390b158d966666f033dbd9f76dc2f4a5ed40/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", [System.Type]::Missing, [System.Type]::Missing, "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/063ede46-8130-4a50-9e03-6494f895b9aa.md", [System.Type]::Missing, [System.Type]::Missing, "063ede46-8130-4a50-9e03-6494f895b9aa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6e783c168c920fc8ee2c92b570627276ee1d6b/e2e/455c1867-16d2-41d4-b97b-ebfc6137f378.md", [System.Type]::Missing, [System.Type]::Missing, "455c1867-16d2-41d4-b97b-ebfc6137f378.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/b9d5e119-8939-4ab9-ac75-292046d3d8d1.md", [System.Type]::Missing, [System.Type]::Missing, "b9d5e119-8939-4ab9-ac75-292046d3d8d1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/897158552f04e39013c8f626cae66f5d9e2f0a5b/e2e/dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md", [System.Type]::Missing, [System.Type]::Missing, "dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d3feedbe1faaf56c3e3a45b046ea53ccb6f815d/e2e/df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md", [System.Type]::Missing, [System.Type]::Missing, "df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f3952b8db7eb72545c9b582a743d2791ac23b108/e2e/df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md", [System.Type]::Missing, [System.Type]::Missing, "df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d482a0e3156d40b34d92c6dcf3f49c866658d47/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md", [System.Type]::Missing, [System.Type]::Missing, "72446337-f112-42e2-8cfc-c5b42674408a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49c8ff281eaf9230f89801459ee5727714e6c121/e2e/72446337-f112-42e2-8cfc-c5b42674408a.md", [System.Type]::Missing, [System.Type]::Missing, "e39f574c-b915-4285-95c4-dfdd9da38f93.md") | Out-Null

# =====================  Sheet "de-de"  ======================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A8").Value = "72446337-f112-42e2-8cfc-c5b42674408a.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "72446337-f112-42e2-8cfc-c5b42674408a.24fd16ff902a084a4561b67b78c898624ad87e36.de-de.xlf"
$ws.Range("H8").Value = "2016-09-07 02:55:41"

$ws.Range("A9").Value = "e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("G9").Value = "e39f574c-b915-4285-95c4-dfdd9da38f93.6f9994cd4a132697e5c945d3f283e9462aaed486.de-de.xlf"
$ws.Range("H9").Value = "2016-09-07 03:08:41"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe85018917a4b35c6d742d73d641a78584ef05a/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", [System.Type]::Missing, [System.Type]::Missing, "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4263691dba76e0635ca6a2ce624d708486279802/e2e/73367972-c8e9-41d0-a2c7-5c5ca94699f6.md", [System.Type]::Missing, [System.Type]::Missing, "73367972-c8e9-41d0-a2c7-5c5ca94699f6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/063ede46-8130-4a50-9e03-6494f895b9aa.md", [System.Type]::Missing, [System.Type]::Missing, "063ede46-8130-4a50-9e03-6494f895b9aa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6e783c168c920fc8ee2c92b570627276ee1d6b/e2e/455c1867-16d2-41d4-b97b-ebfc6137f378.md", [System.Type]::Missing, [System.Type]::Missing, "455c1867-16d2-41d4-b97b-ebfc6137f378.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea58604a2b7a8bb941ffe13a0918c1e5cc390f8a/e2e/b9d5e119-8939-4ab9-ac75-292046d3d8d1.md", [System.Type]::Missing, [System.Type]::Missing, "b9d5e119-8939-4ab9-ac75-292046d3d8d1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/897158552f04e39013c8f626cae66f5d9e2f0a5b/e2e/dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md", [System.Type]::Missing, [System.Type]::Missing, "dda15af0-6a5a-470b-ba2b-bb84fd6d686e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d3feedbe1faaf56c3e3a45b046ea53ccb6f815d/e2e/df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md", [System.Type]::Missing, [System.Type]::Missing, "df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6dbdfb16a4ed8ab6d98955fa38cf07c0427a4ccf/e2e/df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md", [System.Type]::Missing, [System.Type]::Missing, "df70e7cd-d95b-4bc7-bce5-6e3d67de0faa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d482a0e3156d40b34d92c6dcf3f49c866658d47/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md", [System.Type]::Missing, [System.Type]::Missing, "72446337-f112-42e2-8cfc-c5b42674408a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49c8ff281eaf9230f89801459ee5727714e6c121/e2e/72446337-f112-42e2-8cfc-c5b42674408a.md", [System.Type]::Missing, [System.Type]::Missing, "e39f574c-b915-4285-95c4-dfdd9da38f93.md") | Out-Null
